$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume figures; force text storage via NumberFormat
# so numeric-looking values (e.g. "5.124") aren't silently coerced to
# doubles, then reset the style pointer back to Normal so no new cell
# style is introduced.
function Set-TextValue {
    param($cell, [string]$value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "27.857.37"
Set-TextValue $ws.Range("E2") "  +0.81%  "
# Row 3
Set-TextValue $ws.Range("D3") "1.747.21"
Set-TextValue $ws.Range("E3") "  -1.09%  "
# Row 4
Set-TextValue $ws.Range("D4") "1.006"
Set-TextValue $ws.Range("E4") "  +0.34%  "
# Row 5
Set-TextValue $ws.Range("D5") "333.01"
Set-TextValue $ws.Range("E5") "  -0.84%  "
# Row 6
Set-TextValue $ws.Range("E6") "  +0.23%  "
# Row 7
Set-TextValue $ws.Range("D7") "0.3864"
Set-TextValue $ws.Range("E7") "  +0.60%  "
# Row 8
Set-TextValue $ws.Range("D8") "0.3377"
Set-TextValue $ws.Range("E8") "  -1.19%  "
# Row 9
Set-TextValue $ws.Range("D9") "45.80"
Set-TextValue $ws.Range("E9") "  -2.57%  "
# Row 10
Set-TextValue $ws.Range("D10") "1.113"
Set-TextValue $ws.Range("E10") "  -2.18%  "
# Row 11
Set-TextValue $ws.Range("D11") "0.07184"
Set-TextValue $ws.Range("E11") "  -3.01%  "
# Row 12
Set-TextValue $ws.Range("D12") "1.004"
Set-TextValue $ws.Range("E12") "  +0.36%  "
# Row 13
Set-TextValue $ws.Range("D13") "22.43"
Set-TextValue $ws.Range("E13") "  +0.33%  "
# Row 14
Set-TextValue $ws.Range("D14") "6.172"
Set-TextValue $ws.Range("E14") "  -2.85%  "
# Row 15
Set-TextValue $ws.Range("D15") "1.748.94"
Set-TextValue $ws.Range("E15") "  -0.99%  "
# Row 16
Set-TextValue $ws.Range("D16") "7.071"
Set-TextValue $ws.Range("E16") "  -0.18%  "
# Row 17
Set-TextValue $ws.Range("D17") "0.00001057"
Set-TextValue $ws.Range("E17") "  -1.49%  "
# Row 18
Set-TextValue $ws.Range("D18") "0.06602"
Set-TextValue $ws.Range("E18") "  -0.75%  "
# Row 19
Set-TextValue $ws.Range("D19") "79.24"
Set-TextValue $ws.Range("E19") "  -3.52%  "
# Row 20
Set-TextValue $ws.Range("D20") "1.002"
Set-TextValue $ws.Range("E20") "  +0.26%  "
# Row 21
Set-TextValue $ws.Range("D21") "16.74"
Set-TextValue $ws.Range("E21") "  -3.41%  "
# Row 22
Set-TextValue $ws.Range("D22") "6.166"
Set-TextValue $ws.Range("E22") "  -3.68%  "
# Row 23
Set-TextValue $ws.Range("D23") "27.878.92"
Set-TextValue $ws.Range("E23") "  +0.84%  "
# Row 24
Set-TextValue $ws.Range("D24") "11.63"
Set-TextValue $ws.Range("E24") "  -3.35%  "
# Row 25
Set-TextValue $ws.Range("D25") "2.397"
Set-TextValue $ws.Range("E25") "  +0.34%  "
# Row 26
Set-TextValue $ws.Range("D26") "153.89"
Set-TextValue $ws.Range("E26") "  +1.23%  "
# Row 27
Set-TextValue $ws.Range("D27") "19.82"
Set-TextValue $ws.Range("E27") "  -3.97%  "
# Row 28
Set-TextValue $ws.Range("D28") "2.297"
Set-TextValue $ws.Range("E28") "  -4.91%  "
# Row 29
Set-TextValue $ws.Range("D29") "1.949.24"
Set-TextValue $ws.Range("E29") "  -0.87%  "
# Row 30
Set-TextValue $ws.Range("D30") "1.289"
Set-TextValue $ws.Range("E30") "  -9.67%  "
# Row 31
Set-TextValue $ws.Range("D31") "131.03"
Set-TextValue $ws.Range("E31") "  -2.50%  "
# Row 32
Set-TextValue $ws.Range("D32") "4.025"
Set-TextValue $ws.Range("E32") "  +1.62%  "
# Row 33
Set-TextValue $ws.Range("D33") "5.817"
Set-TextValue $ws.Range("E33") "  -4.99%  "
# Row 34
Set-TextValue $ws.Range("D34") "0.08783"
Set-TextValue $ws.Range("E34") "  -0.20%  "
# Row 35
Set-TextValue $ws.Range("D35") "12.09"
Set-TextValue $ws.Range("E35") "  -5.13%  "
# Row 36
Set-TextValue $ws.Range("D36") "1.537"
Set-TextValue $ws.Range("E36") "  +1.87%  "
# Row 37
Set-TextValue $ws.Range("D37") "0.6524"
Set-TextValue $ws.Range("E37") "  -3.90%  "
# Row 38
Set-TextValue $ws.Range("B38") "VeChain"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D38") "0.02275"
Set-TextValue $ws.Range("E38") "  -5.94%  "
# Row 39
Set-TextValue $ws.Range("B39") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D39") "5.124"
Set-TextValue $ws.Range("E39") "  -3.68%  "
# Row 40
Set-TextValue $ws.Range("D40") "0.06121"
Set-TextValue $ws.Range("E40") "  -2.99%  "
# Row 41
Set-TextValue $ws.Range("D41") "0.2092"
Set-TextValue $ws.Range("E41") "  -4.14%  "
# Row 42
Set-TextValue $ws.Range("E42") "  -3.83%  "
# Row 43
Set-TextValue $ws.Range("D43") "8.004"
Set-TextValue $ws.Range("E43") "  -3.16%  "
# Row 44
Set-TextValue $ws.Range("E44") "  +0.23%  "
# Row 45
Set-TextValue $ws.Range("D45") "13.69"
Set-TextValue $ws.Range("E45") "  -3.53%  "
# Row 46
Set-TextValue $ws.Range("D46") "3.816"
Set-TextValue $ws.Range("E46") "  -0.81%  "
# Row 47
Set-TextValue $ws.Range("D47") "0.6028"
# Row 48
Set-TextValue $ws.Range("D48") "126.68"
Set-TextValue $ws.Range("E48") "  -3.80%  "
# Row 49
Set-TextValue $ws.Range("D49") "1.994"
Set-TextValue $ws.Range("E49") "  -3.84%  "
# Row 50
Set-TextValue $ws.Range("D50") "1.163"
Set-TextValue $ws.Range("E50") "  +1.74%  "
# Row 51
Set-TextValue $ws.Range("D51") "1.103"
Set-TextValue $ws.Range("E51") "  +3.84%  "
